$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.801.60'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.54%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.891.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.21%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7753'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3121'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.22'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07165'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08067'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7639'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.05%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.937.52'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.85%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.445'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.20'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.157'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.789.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007758'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9995'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.092'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.109.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.95%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1601'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.379'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.58'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.71'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.18%  '
$ws.Range("E29").Value = '  -3.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.430'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.47%  '
$ws.Range("E31").Value = '  +0.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.466'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.096'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05518'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.261'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7433'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9973'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("E38").Value = '  -3.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01911'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.779'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.137.62'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '73.70'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4417'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.53%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8535'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.32%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.853'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9996'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.886'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.949'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.430'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.014'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +10.11%  '
